$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1187.9
$ws.Cells.Item(2, 3).Value2 = 1208.68
$ws.Cells.Item(2, 4).Value2 = 1167.12
$ws.Cells.Item(2, 5).Value2 = 1199.5
$ws.Cells.Item(2, 6).Value2 = 1176.3
$ws.Cells.Item(3, 2).Value2 = 1340.6
$ws.Cells.Item(3, 3).Value2 = 1367.27
$ws.Cells.Item(3, 4).Value2 = 1313.93
$ws.Cells.Item(3, 5).Value2 = 1354
$ws.Cells.Item(3, 6).Value2 = 1327.2
$ws.Cells.Item(3, 7).Value2 = 1.129
$ws.Cells.Item(4, 2).Value2 = 1040.69
$ws.Cells.Item(4, 3).Value2 = 1058.9
$ws.Cells.Item(4, 4).Value2 = 1022.48
$ws.Cells.Item(4, 5).Value2 = 1050.88
$ws.Cells.Item(4, 6).Value2 = 1030.5
$ws.Cells.Item(4, 7).Value2 = 0.876
$ws.Cells.Item(5, 2).Value2 = 185.88
$ws.Cells.Item(5, 3).Value2 = 195.17
$ws.Cells.Item(5, 4).Value2 = 176.59
$ws.Cells.Item(5, 5).Value2 = 187.73
$ws.Cells.Item(5, 6).Value2 = 184.03
$ws.Cells.Item(5, 7).Value2 = 0.157
$ws.Cells.Item(6, 2).Value2 = 152.48
$ws.Cells.Item(6, 3).Value2 = 155.48
$ws.Cells.Item(6, 4).Value2 = 149.48
$ws.Cells.Item(6, 5).Value2 = 154
$ws.Cells.Item(6, 6).Value2 = 150.96
$ws.Cells.Item(6, 7).Value2 = 0.128
$ws.Cells.Item(7, 2).Value2 = 42.72
$ws.Cells.Item(7, 3).Value2 = 48.31
$ws.Cells.Item(7, 4).Value2 = 39.73
$ws.Cells.Item(7, 7).Value2 = 0.036
$ws.Cells.Item(8, 2).Value2 = 1596.83
$ws.Cells.Item(8, 3).Value2 = 1628.28
$ws.Cells.Item(8, 4).Value2 = 1565.38
$ws.Cells.Item(8, 5).Value2 = 1612.79
$ws.Cells.Item(8, 6).Value2 = 1580.87
$ws.Cells.Item(8, 7).Value2 = 1.344
$ws.Cells.Item(9, 2).Value2 = 3085.53
$ws.Cells.Item(9, 3).Value2 = 3360.14
$ws.Cells.Item(9, 4).Value2 = 2900.4
$ws.Cells.Item(9, 7).Value2 = 2.598
$ws.Cells.Item(10, 2).Value2 = 939.5
$ws.Cells.Item(10, 3).Value2 = 958
$ws.Cells.Item(10, 4).Value2 = 921
$ws.Cells.Item(10, 5).Value2 = 948.89
$ws.Cells.Item(10, 6).Value2 = 930.11
$ws.Cells.Item(10, 7).Value2 = 0.791
$ws.Cells.Item(11, 2).Value2 = 1278.76
$ws.Cells.Item(11, 3).Value2 = 1303.95
$ws.Cells.Item(11, 4).Value2 = 1253.57
$ws.Cells.Item(11, 5).Value2 = 1291.54
$ws.Cells.Item(11, 6).Value2 = 1265.98
$ws.Cells.Item(11, 7).Value2 = 1.077
$ws.Cells.Item(12, 2).Value2 = 132.76
$ws.Cells.Item(12, 3).Value2 = 136.01
$ws.Cells.Item(12, 4).Value2 = 129.51
$ws.Cells.Item(12, 5).Value2 = 134.08
$ws.Cells.Item(12, 6).Value2 = 131.44
$ws.Cells.Item(12, 7).Value2 = 0.112
$ws.Cells.Item(13, 2).Value2 = 859.6799999999999
$ws.Cells.Item(13, 3).Value2 = 876.61
$ws.Cells.Item(13, 4).Value2 = 842.75
$ws.Cells.Item(13, 5).Value2 = 868.27
$ws.Cells.Item(13, 6).Value2 = 851.09
$ws.Cells.Item(13, 7).Value2 = 0.724
$ws.Cells.Item(14, 2).Value2 = 831.17
$ws.Cells.Item(14, 3).Value2 = 847.54
$ws.Cells.Item(14, 4).Value2 = 814.8
$ws.Cells.Item(14, 5).Value2 = 839.48
$ws.Cells.Item(14, 6).Value2 = 822.86
$ws.Cells.Item(14, 7).Value2 = 0.7
$ws.Cells.Item(15, 2).Value2 = 52.75
$ws.Cells.Item(15, 3).Value2 = 57.23
$ws.Cells.Item(15, 4).Value2 = 48.01
$ws.Cells.Item(15, 5).Value2 = 53.33
$ws.Cells.Item(15, 6).Value2 = 52.17
$ws.Cells.Item(15, 7).Value2 = 0.044
$ws.Cells.Item(16, 2).Value2 = 1.43
$ws.Cells.Item(16, 3).Value2 = 1.54
$ws.Cells.Item(16, 4).Value2 = 1.32
$ws.Cells.Item(17, 2).Value2 = 106.16
$ws.Cells.Item(17, 5).Value2 = 107.32
$ws.Cells.Item(17, 6).Value2 = 105
$ws.Cells.Item(17, 7).Value2 = 0.089
$ws.Cells.Item(18, 2).Value2 = 0.42
$ws.Cells.Item(19, 2).Value2 = 383.64
$ws.Cells.Item(19, 3).Value2 = 422
$ws.Cells.Item(19, 4).Value2 = 352.95
$ws.Cells.Item(19, 7).Value2 = 0.323
$ws.Cells.Item(20, 2).Value2 = 180.27
$ws.Cells.Item(20, 3).Value2 = 184.68
$ws.Cells.Item(20, 4).Value2 = 175.86
$ws.Cells.Item(20, 5).Value2 = 182.07
$ws.Cells.Item(20, 6).Value2 = 178.47
$ws.Cells.Item(20, 7).Value2 = 0.152
$ws.Cells.Item(21, 2).Value2 = 133.24
$ws.Cells.Item(21, 3).Value2 = 136.5
$ws.Cells.Item(21, 4).Value2 = 129.98
$ws.Cells.Item(21, 5).Value2 = 134.57
$ws.Cells.Item(21, 6).Value2 = 131.91
$ws.Cells.Item(21, 7).Value2 = 0.112
$ws.Cells.Item(22, 2).Value2 = 316.68
$ws.Cells.Item(22, 3).Value2 = 336.63
$ws.Cells.Item(22, 4).Value2 = 294.83
$ws.Cells.Item(22, 5).Value2 = 319.84
$ws.Cells.Item(22, 6).Value2 = 313.52
$ws.Cells.Item(23, 2).Value2 = 3926.29
$ws.Cells.Item(23, 3).Value2 = 4181.49
$ws.Cells.Item(23, 4).Value2 = 3612.19
$ws.Cells.Item(23, 5).Value2 = 3965.55
$ws.Cells.Item(23, 6).Value2 = 3887.03
$ws.Cells.Item(23, 7).Value2 = 3.305
$ws.Cells.Item(24, 2).Value2 = 3150.93
$ws.Cells.Item(24, 3).Value2 = 3352.58
$ws.Cells.Item(24, 4).Value2 = 2898.86
$ws.Cells.Item(24, 5).Value2 = 3182.43
$ws.Cells.Item(24, 6).Value2 = 3119.43
$ws.Cells.Item(24, 7).Value2 = 2.653
$ws.Cells.Item(25, 2).Value2 = 323.41
$ws.Cells.Item(25, 3).Value2 = 341.19
$ws.Cells.Item(25, 4).Value2 = 301.1
$ws.Cells.Item(25, 5).Value2 = 326.64
$ws.Cells.Item(25, 6).Value2 = 320.18
$ws.Cells.Item(26, 2).Value2 = 1677.23
$ws.Cells.Item(26, 3).Value2 = 1826.5
$ws.Cells.Item(26, 4).Value2 = 1543.06
$ws.Cells.Item(26, 7).Value2 = 1.412
$ws.Cells.Item(27, 2).Value2 = 75.66
$ws.Cells.Item(27, 7).Value2 = 0.064
$ws.Cells.Item(28, 2).Value2 = 36.17
$ws.Cells.Item(28, 3).Value2 = 37.97
$ws.Cells.Item(28, 4).Value2 = 34
$ws.Cells.Item(28, 5).Value2 = 36.53
$ws.Cells.Item(28, 6).Value2 = 35.81
$ws.Cells.Item(28, 7).Value2 = 0.03
$ws.Cells.Item(29, 2).Value2 = 871.95
$ws.Cells.Item(29, 3).Value2 = 889.3
$ws.Cells.Item(29, 4).Value2 = 854.6
$ws.Cells.Item(29, 5).Value2 = 880.66
$ws.Cells.Item(29, 6).Value2 = 863.24
$ws.Cells.Item(29, 7).Value2 = 0.734
$ws.Cells.Item(30, 2).Value2 = 283.64
$ws.Cells.Item(30, 3).Value2 = 301.5
$ws.Cells.Item(30, 4).Value2 = 262.66
$ws.Cells.Item(30, 6).Value2 = 280.81
$ws.Cells.Item(30, 7).Value2 = 0.239
$ws.Cells.Item(31, 2).Value2 = 8.359999999999999
$ws.Cells.Item(31, 3).Value2 = 8.94
$ws.Cells.Item(31, 4).Value2 = 7.53
$ws.Cells.Item(31, 5).Value2 = 8.44
$ws.Cells.Item(31, 6).Value2 = 8.279999999999999
$ws.Cells.Item(31, 7).Value2 = 0.007
$ws.Cells.Item(32, 2).Value2 = 324.3
$ws.Cells.Item(32, 7).Value2 = 0.273
$ws.Cells.Item(33, 2).Value2 = 2.74
$ws.Cells.Item(34, 2).Value2 = 871.95
$ws.Cells.Item(34, 3).Value2 = 906.8200000000001
$ws.Cells.Item(34, 4).Value2 = 819.64
$ws.Cells.Item(34, 7).Value2 = 0.734
$ws.Cells.Item(35, 2).Value2 = 15.98
$ws.Cells.Item(36, 2).Value2 = 6.79
$ws.Cells.Item(37, 2).Value2 = 13.85
$ws.Cells.Item(38, 2).Value2 = 23.41
$ws.Cells.Item(38, 3).Value2 = 25.75
$ws.Cells.Item(38, 4).Value2 = 21.5
$ws.Cells.Item(38, 5).Value2 = 23.64
$ws.Cells.Item(38, 6).Value2 = 23.18
$ws.Cells.Item(38, 7).Value2 = 0.02
$ws.Cells.Item(39, 2).Value2 = 56.98
$ws.Cells.Item(39, 3).Value2 = 62.67
$ws.Cells.Item(39, 4).Value2 = 52.2
$ws.Cells.Item(39, 5).Value2 = 57.54
$ws.Cells.Item(39, 6).Value2 = 56.42
$ws.Cells.Item(39, 7).Value2 = 0.048
$ws.Cells.Item(40, 2).Value2 = 211.61
$ws.Cells.Item(40, 3).Value2 = 233.19
$ws.Cells.Item(40, 4).Value2 = 194.69
$ws.Cells.Item(40, 6).Value2 = 209.08
$ws.Cells.Item(40, 7).Value2 = 0.178
$ws.Cells.Item(41, 2).Value2 = 5.24
$ws.Cells.Item(41, 3).Value2 = 5.85
$ws.Cells.Item(41, 4).Value2 = 4.63
$ws.Cells.Item(41, 5).Value2 = 5.29
$ws.Cells.Item(41, 6).Value2 = 5.19
$ws.Cells.Item(42, 2).Value2 = 75.45
$ws.Cells.Item(42, 3).Value2 = 79.97
$ws.Cells.Item(42, 4).Value2 = 69.42
$ws.Cells.Item(42, 5).Value2 = 76.34999999999999
$ws.Cells.Item(42, 6).Value2 = 74.55
$ws.Cells.Item(42, 7).Value2 = 0.064
$ws.Cells.Item(43, 2).Value2 = 16.16
$ws.Cells.Item(43, 3).Value2 = 17.29
$ws.Cells.Item(43, 4).Value2 = 14.39
$ws.Cells.Item(43, 5).Value2 = 16.32
$ws.Cells.Item(43, 6).Value2 = 16
$ws.Cells.Item(43, 7).Value2 = 0.014
$ws.Cells.Item(44, 2).Value2 = 3.65
$ws.Cells.Item(44, 3).Value2 = 3.98
$ws.Cells.Item(44, 4).Value2 = 3.36
$ws.Cells.Item(44, 5).Value2 = 3.69
$ws.Cells.Item(44, 6).Value2 = 3.61
$ws.Cells.Item(45, 2).Value2 = 285.72
$ws.Cells.Item(45, 3).Value2 = 308.57
$ws.Cells.Item(45, 4).Value2 = 262.87
$ws.Cells.Item(45, 5).Value2 = 288.86
$ws.Cells.Item(45, 6).Value2 = 282.58
$ws.Cells.Item(45, 7).Value2 = 0.241
